# Update NATMI LR-pair (Hras-Agtr1a) output sheet with recomputed TPM-derived
# ligand/receptor expression, specificity and edge-weight statistics.
# All 20 columns (A:T) keep their row/column layout; only the numeric result
# columns G,H,I,J,M,N,O,P,Q,R,S,T for rows 2-16 are refreshed with values
# produced by the updated "new tpm" scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 9.697207333333333
$ws.Cells.Item(2, 8).Value = 29.091622
$ws.Cells.Item(2, 9).Value = 0.3243108558382838
$ws.Cells.Item(2, 10).Value = 0.3243108558382838
$ws.Cells.Item(2, 13).Value = 0.6100786666666667
$ws.Cells.Item(2, 14).Value = 1.830236
$ws.Cells.Item(2, 15).Value = 0.06402955811028149
$ws.Cells.Item(2, 16).Value = 0.06402955811028149
$ws.Cells.Item(2, 17).Value = 5.916059320310222
$ws.Cells.Item(2, 18).Value = 53.244533882792
$ws.Cells.Item(2, 19).Value = 0.02076548078969252
$ws.Cells.Item(2, 20).Value = 0.02076548078969252
$ws.Cells.Item(3, 7).Value = 9.697207333333333
$ws.Cells.Item(3, 8).Value = 29.091622
$ws.Cells.Item(3, 9).Value = 0.3243108558382838
$ws.Cells.Item(3, 10).Value = 0.3243108558382838
$ws.Cells.Item(3, 13).Value = 7.236132333333333
$ws.Cells.Item(3, 15).Value = 0.7594534623909487
$ws.Cells.Item(3, 16).Value = 0.7594534623909487
$ws.Cells.Item(3, 17).Value = 70.17027552777044
$ws.Cells.Item(3, 18).Value = 631.5324797499339
$ws.Cells.Item(3, 19).Value = 0.2462990023573565
$ws.Cells.Item(3, 20).Value = 0.2462990023573564
$ws.Cells.Item(4, 7).Value = 9.697207333333333
$ws.Cells.Item(4, 8).Value = 29.091622
$ws.Cells.Item(4, 9).Value = 0.3243108558382838
$ws.Cells.Item(4, 10).Value = 0.3243108558382838
$ws.Cells.Item(4, 13).Value = 1.681867666666667
$ws.Cells.Item(4, 14).Value = 5.045603
$ws.Cells.Item(4, 15).Value = 0.1765169794987699
$ws.Cells.Item(4, 16).Value = 0.1765169794987699
$ws.Cells.Item(4, 17).Value = 16.30941947089622
$ws.Cells.Item(4, 18).Value = 146.784775238066
$ws.Cells.Item(4, 19).Value = 0.05724637269123486
$ws.Cells.Item(4, 20).Value = 0.05724637269123486
$ws.Cells.Item(5, 9).Value = 0.2826325233457075
$ws.Cells.Item(5, 10).Value = 0.2826325233457074
$ws.Cells.Item(5, 13).Value = 0.6100786666666667
$ws.Cells.Item(5, 14).Value = 1.830236
$ws.Cells.Item(5, 15).Value = 0.06402955811028149
$ws.Cells.Item(5, 16).Value = 0.06402955811028149
$ws.Cells.Item(5, 17).Value = 5.15576566082
$ws.Cells.Item(5, 18).Value = 46.40189094738
$ws.Cells.Item(5, 19).Value = 0.01809683557741947
$ws.Cells.Item(5, 20).Value = 0.01809683557741947
$ws.Cells.Item(6, 9).Value = 0.2826325233457075
$ws.Cells.Item(6, 10).Value = 0.2826325233457074
$ws.Cells.Item(6, 13).Value = 7.236132333333333
$ws.Cells.Item(6, 15).Value = 0.7594534623909487
$ws.Cells.Item(6, 16).Value = 0.7594534623909487
$ws.Cells.Item(6, 18).Value = 550.3720122631349
$ws.Cells.Item(6, 19).Value = 0.2146462484391882
$ws.Cells.Item(6, 20).Value = 0.2146462484391881
$ws.Cells.Item(7, 9).Value = 0.2826325233457075
$ws.Cells.Item(7, 10).Value = 0.2826325233457074
$ws.Cells.Item(7, 13).Value = 1.681867666666667
$ws.Cells.Item(7, 14).Value = 5.045603
$ws.Cells.Item(7, 15).Value = 0.1765169794987699
$ws.Cells.Item(7, 16).Value = 0.1765169794987699
$ws.Cells.Item(7, 17).Value = 14.213438422985
$ws.Cells.Item(7, 18).Value = 127.920945806865
$ws.Cells.Item(7, 19).Value = 0.04988943932909986
$ws.Cells.Item(7, 20).Value = 0.04988943932909985
$ws.Cells.Item(8, 7).Value = 3.910524
$ws.Cells.Item(8, 8).Value = 11.731572
$ws.Cells.Item(8, 9).Value = 0.1307825378608469
$ws.Cells.Item(8, 10).Value = 0.1307825378608469
$ws.Cells.Item(8, 13).Value = 0.6100786666666667
$ws.Cells.Item(8, 14).Value = 1.830236
$ws.Cells.Item(8, 15).Value = 0.06402955811028149
$ws.Cells.Item(8, 16).Value = 0.06402955811028149
$ws.Cells.Item(8, 17).Value = 2.385727267888
$ws.Cells.Item(8, 18).Value = 21.471545410992
$ws.Cells.Item(8, 19).Value = 0.008373948107771187
$ws.Cells.Item(8, 20).Value = 0.008373948107771185
$ws.Cells.Item(9, 7).Value = 3.910524
$ws.Cells.Item(9, 8).Value = 11.731572
$ws.Cells.Item(9, 9).Value = 0.1307825378608469
$ws.Cells.Item(9, 10).Value = 0.1307825378608469
$ws.Cells.Item(9, 13).Value = 7.236132333333333
$ws.Cells.Item(9, 15).Value = 0.7594534623909487
$ws.Cells.Item(9, 16).Value = 0.7594534623909487
$ws.Cells.Item(9, 17).Value = 28.297069156676
$ws.Cells.Item(9, 18).Value = 254.673622410084
$ws.Cells.Item(9, 19).Value = 0.09932325119869553
$ws.Cells.Item(9, 20).Value = 0.09932325119869551
$ws.Cells.Item(10, 7).Value = 3.910524
$ws.Cells.Item(10, 8).Value = 11.731572
$ws.Cells.Item(10, 9).Value = 0.1307825378608469
$ws.Cells.Item(10, 10).Value = 0.1307825378608469
$ws.Cells.Item(10, 13).Value = 1.681867666666667
$ws.Cells.Item(10, 14).Value = 5.045603
$ws.Cells.Item(10, 15).Value = 0.1765169794987699
$ws.Cells.Item(10, 16).Value = 0.1765169794987699
$ws.Cells.Item(10, 17).Value = 6.576983875323999
$ws.Cells.Item(10, 18).Value = 59.192854877916
$ws.Cells.Item(10, 19).Value = 0.02308533855438021
$ws.Cells.Item(10, 20).Value = 0.02308533855438021
$ws.Cells.Item(11, 7).Value = 3.066674
$ws.Cells.Item(11, 8).Value = 9.200022000000001
$ws.Cells.Item(11, 9).Value = 0.1025610400324547
$ws.Cells.Item(11, 10).Value = 0.1025610400324547
$ws.Cells.Item(11, 13).Value = 0.6100786666666667
$ws.Cells.Item(11, 14).Value = 1.830236
$ws.Cells.Item(11, 15).Value = 0.06402955811028149
$ws.Cells.Item(11, 16).Value = 0.06402955811028149
$ws.Cells.Item(11, 17).Value = 1.870912385021334
$ws.Cells.Item(11, 18).Value = 16.838211465192
$ws.Cells.Item(11, 19).Value = 0.006566938072608965
$ws.Cells.Item(11, 20).Value = 0.006566938072608964
$ws.Cells.Item(12, 7).Value = 3.066674
$ws.Cells.Item(12, 8).Value = 9.200022000000001
$ws.Cells.Item(12, 9).Value = 0.1025610400324547
$ws.Cells.Item(12, 10).Value = 0.1025610400324547
$ws.Cells.Item(12, 13).Value = 7.236132333333333
$ws.Cells.Item(12, 15).Value = 0.7594534623909487
$ws.Cells.Item(12, 16).Value = 0.7594534623909487
$ws.Cells.Item(12, 17).Value = 22.19085888719267
$ws.Cells.Item(12, 18).Value = 199.717729984734
$ws.Cells.Item(12, 19).Value = 0.07789033695906443
$ws.Cells.Item(12, 20).Value = 0.0778903369590644
$ws.Cells.Item(13, 7).Value = 3.066674
$ws.Cells.Item(13, 8).Value = 9.200022000000001
$ws.Cells.Item(13, 9).Value = 0.1025610400324547
$ws.Cells.Item(13, 10).Value = 0.1025610400324547
$ws.Cells.Item(13, 13).Value = 1.681867666666667
$ws.Cells.Item(13, 14).Value = 5.045603
$ws.Cells.Item(13, 15).Value = 0.1765169794987699
$ws.Cells.Item(13, 16).Value = 0.1765169794987699
$ws.Cells.Item(13, 17).Value = 5.157739844807334
$ws.Cells.Item(13, 18).Value = 46.419658603266
$ws.Cells.Item(13, 19).Value = 0.01810376500078132
$ws.Cells.Item(13, 20).Value = 0.01810376500078132
$ws.Cells.Item(14, 7).Value = 4.775574
$ws.Cells.Item(14, 8).Value = 14.326722
$ws.Cells.Item(14, 9).Value = 0.1597130429227071
$ws.Cells.Item(14, 10).Value = 0.159713042922707
$ws.Cells.Item(14, 13).Value = 0.6100786666666667
$ws.Cells.Item(14, 14).Value = 1.830236
$ws.Cells.Item(14, 15).Value = 0.06402955811028149
$ws.Cells.Item(14, 16).Value = 0.06402955811028149
$ws.Cells.Item(14, 17).Value = 2.913475818488
$ws.Cells.Item(14, 18).Value = 26.221282366392
$ws.Cells.Item(14, 19).Value = 0.01022635556278935
$ws.Cells.Item(14, 20).Value = 0.01022635556278935
$ws.Cells.Item(15, 7).Value = 4.775574
$ws.Cells.Item(15, 8).Value = 14.326722
$ws.Cells.Item(15, 9).Value = 0.1597130429227071
$ws.Cells.Item(15, 10).Value = 0.159713042922707
$ws.Cells.Item(15, 13).Value = 7.236132333333333
$ws.Cells.Item(15, 15).Value = 0.7594534623909487
$ws.Cells.Item(15, 16).Value = 0.7594534623909487
$ws.Cells.Item(15, 17).Value = 34.55668543162599
$ws.Cells.Item(15, 18).Value = 311.010168884634
$ws.Cells.Item(15, 19).Value = 0.1212946234366441
$ws.Cells.Item(15, 20).Value = 0.1212946234366441
$ws.Cells.Item(16, 7).Value = 4.775574
$ws.Cells.Item(16, 8).Value = 14.326722
$ws.Cells.Item(16, 9).Value = 0.1597130429227071
$ws.Cells.Item(16, 10).Value = 0.159713042922707
$ws.Cells.Item(16, 13).Value = 1.681867666666667
$ws.Cells.Item(16, 14).Value = 5.045603
$ws.Cells.Item(16, 15).Value = 0.1765169794987699
$ws.Cells.Item(16, 16).Value = 0.1765169794987699
$ws.Cells.Item(16, 17).Value = 8.031883500373999
$ws.Cells.Item(16, 18).Value = 72.286951503366
$ws.Cells.Item(16, 19).Value = 0.02819206392327364
$ws.Cells.Item(16, 20).Value = 0.02819206392327364

Write-Host "Updated 168 cells (G2:T16 result columns) with refreshed TPM values."
